$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (wave 0001)
$ws.Range("D6").Value = "'3"
$ws.Range("E6").Value = "'25"
$ws.Range("G6").Value = "'12"

# Row 7 (wave 0002)
$ws.Range("D7").Value = "'35"
$ws.Range("E7").Value = "'35"
$ws.Range("G7").Value = "'10"

# Row 8 (wave 0003)
$ws.Range("D8").Value = "'80"
$ws.Range("E8").Value = "'30"

# Row 9 (wave 0004)
$ws.Range("D9").Value = "'120"
$ws.Range("E9").Value = "'40"

# Row 10 (wave 0005)
$ws.Range("D10").Value = "'170"
$ws.Range("E10").Value = "'35"
$ws.Range("G10").Value = "'14"

# Row 11 (wave 0006)
$ws.Range("D11").Value = "'220"
$ws.Range("E11").Value = "'45"
